# Update the LR-pair TPM-derived metrics on the active sheet to reflect
# the newly recomputed TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.7198826633506309
$ws.Range("J2").Value = 0.7198826633506309
$ws.Range("M2").Value = 10.92359866666667
$ws.Range("N2").Value = 32.770796
$ws.Range("O2").Value = 0.2236009040380497
$ws.Range("P2").Value = 0.2236009040380497
$ws.Range("Q2").Value = 0.6576734637244446
$ws.Range("R2").Value = 5.919061173520001
$ws.Range("S2").Value = 0.1609664143265201
$ws.Range("T2").Value = 0.16096641432652

# Row 3
$ws.Range("I3").Value = 0.7198826633506309
$ws.Range("J3").Value = 0.7198826633506309
$ws.Range("O3").Value = 0.4261214970992155
$ws.Range("P3").Value = 0.4261214970992155
$ws.Range("S3").Value = 0.3067574782427414
$ws.Range("T3").Value = 0.3067574782427414

# Row 4
$ws.Range("I4").Value = 0.7198826633506309
$ws.Range("J4").Value = 0.7198826633506309
$ws.Range("M4").Value = 13.06524766666667
$ws.Range("N4").Value = 39.195743
$ws.Range("O4").Value = 0.2674394472823625
$ws.Range("P4").Value = 0.2674394472823625
$ws.Range("Q4").Value = 0.7866150111844444
$ws.Range("R4").Value = 7.07953510066
$ws.Range("S4").Value = 0.1925250215946478
$ws.Range("T4").Value = 0.1925250215946478

# Row 5
$ws.Range("I5").Value = 0.7198826633506309
$ws.Range("J5").Value = 0.7198826633506309
$ws.Range("M5").Value = 4.046901
$ws.Range("N5").Value = 12.140703
$ws.Range("O5").Value = 0.0828381515803724
$ws.Range("P5").Value = 0.0828381515803724
$ws.Range("Q5").Value = 0.24365041954
$ws.Range("R5").Value = 2.19285377586
$ws.Range("S5").Value = 0.05963374918672176
$ws.Range("T5").Value = 0.05963374918672176

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02342733333333333
$ws.Range("H6").Value = 0.070282
$ws.Range("I6").Value = 0.2801173366493691
$ws.Range("J6").Value = 0.2801173366493691
$ws.Range("M6").Value = 10.92359866666667
$ws.Range("N6").Value = 32.770796
$ws.Range("O6").Value = 0.2236009040380497
$ws.Range("P6").Value = 0.2236009040380497
$ws.Range("Q6").Value = 0.2559107871635556
$ws.Range("R6").Value = 2.303197084472
$ws.Range("S6").Value = 0.06263448971152963
$ws.Range("T6").Value = 0.06263448971152963

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.02342733333333333
$ws.Range("H7").Value = 0.070282
$ws.Range("I7").Value = 0.2801173366493691
$ws.Range("J7").Value = 0.2801173366493691
$ws.Range("O7").Value = 0.4261214970992155
$ws.Range("P7").Value = 0.4261214970992155
$ws.Range("Q7").Value = 0.4876952006035555
$ws.Range("R7").Value = 4.389256805432
$ws.Range("S7").Value = 0.1193640188564741
$ws.Range("T7").Value = 0.1193640188564741

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.02342733333333333
$ws.Range("H8").Value = 0.070282
$ws.Range("I8").Value = 0.2801173366493691
$ws.Range("J8").Value = 0.2801173366493691
$ws.Range("M8").Value = 13.06524766666667
$ws.Range("N8").Value = 39.195743
$ws.Range("O8").Value = 0.2674394472823625
$ws.Range("P8").Value = 0.2674394472823625
$ws.Range("Q8").Value = 0.3060839121695555
$ws.Range("R8").Value = 2.754755209526
$ws.Range("S8").Value = 0.07491442568771473
$ws.Range("T8").Value = 0.07491442568771473

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.02342733333333333
$ws.Range("H9").Value = 0.070282
$ws.Range("I9").Value = 0.2801173366493691
$ws.Range("J9").Value = 0.2801173366493691
$ws.Range("M9").Value = 4.046901
$ws.Range("N9").Value = 12.140703
$ws.Range("O9").Value = 0.0828381515803724
$ws.Range("P9").Value = 0.0828381515803724
$ws.Range("Q9").Value = 0.09480809869399999
$ws.Range("R9").Value = 0.853272888246
$ws.Range("S9").Value = 0.02320440239365064
$ws.Range("T9").Value = 0.02320440239365064
